# monsters.xlsx — "Broke the import for items into two sets" edit.
#
# Continues the spell_evasion/artifact_annulment/affix_resistance/
# entrancing_chance (AE:AH) ramp past its old cap of 1.0 for rows 319-357
# (the "Hell" difficulty band), picking the +0.07/row progression back up
# from row 318's value of 0.95 (AD keeps the old 1.0 cap — only AE:AH are
# touched). Also moves the saved viewport/selection to where that band was
# being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# --- AE319:AH357: un-cap spell_evasion / artifact_annulment /
#     affix_resistance / entrancing_chance, continuing the +0.07 ramp.
#     Each row holds one repeated value across all four columns, so a
#     single broadcast assignment per row reproduces that.
$ws.Range("AE319:AH319").Value = 1.02
$ws.Range("AE320:AH320").Value = 1.0900000000000001
$ws.Range("AE321:AH321").Value = 1.1599999999999999
$ws.Range("AE322:AH322").Value = 1.23
$ws.Range("AE323:AH323").Value = 1.3
$ws.Range("AE324:AH324").Value = 1.37
$ws.Range("AE325:AH325").Value = 1.44
$ws.Range("AE326:AH326").Value = 1.51
$ws.Range("AE327:AH327").Value = 1.58
$ws.Range("AE328:AH328").Value = 1.65
$ws.Range("AE329:AH329").Value = 1.72
$ws.Range("AE330:AH330").Value = 1.79
$ws.Range("AE331:AH331").Value = 1.86
$ws.Range("AE332:AH332").Value = 1.93
$ws.Range("AE333:AH333").Value = 2
$ws.Range("AE334:AH334").Value = 2.0699999999999998
$ws.Range("AE335:AH335").Value = 2.14
$ws.Range("AE336:AH336").Value = 2.21
$ws.Range("AE337:AH337").Value = 2.2799999999999998
$ws.Range("AE338:AH338").Value = 2.35
$ws.Range("AE339:AH339").Value = 2.42
$ws.Range("AE340:AH340").Value = 2.4900000000000002
$ws.Range("AE341:AH341").Value = 2.56
$ws.Range("AE342:AH342").Value = 2.63
$ws.Range("AE343:AH343").Value = 2.7
$ws.Range("AE344:AH344").Value = 2.77
$ws.Range("AE345:AH345").Value = 2.84
$ws.Range("AE346:AH346").Value = 2.91
$ws.Range("AE347:AH347").Value = 2.98
$ws.Range("AE348:AH348").Value = 3.05
$ws.Range("AE349:AH349").Value = 3.12
$ws.Range("AE350:AH350").Value = 3.19
$ws.Range("AE351:AH351").Value = 3.26
$ws.Range("AE352:AH352").Value = 3.33
$ws.Range("AE353:AH353").Value = 3.4
$ws.Range("AE354:AH354").Value = 3.47
$ws.Range("AE355:AH355").Value = 3.54
$ws.Range("AE356:AH356").Value = 3.61
$ws.Range("AE357:AH357").Value = 3.68

# --- Scroll the saved view toward the band just edited and leave the
#     selection on AE316:AH357 (anchored at AE316), matching where the
#     author's cursor ended up after this pass.
$excel.ActiveWindow.ScrollRow = 319
$excel.ActiveWindow.ScrollColumn = 26
$ws.Range("AE316:AH357").Select()
